# Update quarterly recurrence metrics for row 29 (2025Q3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 48
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = 40
$ws.Range("F29").Value = 1.376936316695353
